$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the two "tasa" lines inside the daily conversion note ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.82 = 14664.12 pesos`n✅ 14664.12 pesos = 3.8 = 942.15 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $conversionText

# --- tasas!N10/O10/N12/O12: refreshed rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 262
$ws2.Range("O10").Value = 3842
$ws2.Range("N12").Value = 3860
$ws2.Range("O12").Value = 248
